$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text change
$ws.Range("D1").Value = '${expected_error_msg}'

# Row 2: shift B2/C2 (both blank quotes) out - D2 stays same text but now fed from A2/D2 only
$ws.Range("B2").Value = $null
$ws.Range("C2").Value = $null

# Row 3: shift left - old B3 (blank) removed, C3 becomes "secret_sauce", D3 unchanged text
$ws.Range("B3").Value = $null
$ws.Range("C3").Value = "secret_sauce"

# Row 4: shift left - old C4 (blank) removed, D4 unchanged value moves from D4 to D4 (same) but B4 stays, so just clear C4
$ws.Range("C4").Value = $null

# Row 5,6,7 unaffected in content (same text), no action needed

# Selection change
$ws.Range("D6").Select()
